$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(10, 8).Value = 10666.667  # ALC!H10
$ws.Cells.Item(10, 10).Value = 9000  # ALC!J10
$ws.Cells.Item(10, 12).Value = 9000  # ALC!L10
$ws.Cells.Item(10, 14).Value = -9586  # ALC!N10

$ws.Cells.Item(13, 8).Value = 1933.3334  # ALC!H13
$ws.Cells.Item(13, 9).Value = 0  # ALC!I13
$ws.Cells.Item(13, 10).Value = 1933.3334  # ALC!J13
$ws.Cells.Item(13, 11).Value = 0  # ALC!K13
$ws.Cells.Item(13, 12).Value = 1933.3334  # ALC!L13
$ws.Cells.Item(13, 13).ClearContents()  # ALC!M13
$ws.Cells.Item(13, 14).Value = -2271.3334  # ALC!N13

$ws.Cells.Item(18, 8).Value = 640.4286  # ALC!H18
$ws.Cells.Item(18, 9).Value = 546.8333  # ALC!I18
$ws.Cells.Item(18, 10).Value = 1202  # ALC!J18
$ws.Cells.Item(18, 11).Value = 546.8333  # ALC!K18
$ws.Cells.Item(18, 12).Value = 1202  # ALC!L18
$ws.Cells.Item(18, 13).Value = -262.8333  # ALC!M18
$ws.Cells.Item(18, 14).Value = -1770  # ALC!N18

$ws.Cells.Item(29, 8).Value = 512.875  # ALC!H29
$ws.Cells.Item(29, 9).Value = 84.333336  # ALC!I29
$ws.Cells.Item(29, 10).Value = 770  # ALC!J29
$ws.Cells.Item(29, 11).Value = 253.000008  # ALC!K29
$ws.Cells.Item(29, 12).Value = 2310  # ALC!L29
$ws.Cells.Item(29, 13).Value = 27.99999199999999  # ALC!M29
$ws.Cells.Item(29, 14).Value = -2872  # ALC!N29

$ws.Cells.Item(31, 8).Value = 3066.818  # ALC!H31
$ws.Cells.Item(31, 9).Value = 575.75  # ALC!I31
$ws.Cells.Item(31, 11).Value = 1727.25  # ALC!K31
$ws.Cells.Item(31, 13).Value = -1497.25  # ALC!M31

$ws.Cells.Item(43, 8).Value = 767.375  # ALC!H43
$ws.Cells.Item(43, 9).Value = 0  # ALC!I43
$ws.Cells.Item(43, 10).Value = 767.375  # ALC!J43
$ws.Cells.Item(43, 11).Value = 0  # ALC!K43
$ws.Cells.Item(43, 12).Value = 767.375  # ALC!L43
$ws.Cells.Item(43, 13).ClearContents()  # ALC!M43
$ws.Cells.Item(43, 14).Value = -905.375  # ALC!N43

$ws.Cells.Item(137, 8).Value = 1426.6666  # ALC!H137
$ws.Cells.Item(137, 9).Value = 1239.8667  # ALC!I137
$ws.Cells.Item(137, 10).Value = 1613.4667  # ALC!J137
$ws.Cells.Item(137, 11).Value = 3719.6001  # ALC!K137
$ws.Cells.Item(137, 12).Value = 4840.4001  # ALC!L137
$ws.Cells.Item(137, 13).Value = -1169.6001  # ALC!M137
$ws.Cells.Item(137, 14).Value = -9940.400099999999  # ALC!N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 12604  # ARM!H26
$ws.Cells.Item(26, 9).Value = 12604  # ARM!I26
$ws.Cells.Item(26, 11).Value = 12604  # ARM!K26
$ws.Cells.Item(26, 13).Value = -12274  # ARM!M26

$ws.Cells.Item(38, 8).Value = 11346.333  # ARM!H38
$ws.Cells.Item(38, 9).Value = 2019.5  # ARM!I38
$ws.Cells.Item(38, 10).Value = 30000  # ARM!J38
$ws.Cells.Item(38, 11).Value = 2019.5  # ARM!K38
$ws.Cells.Item(38, 12).Value = 30000  # ARM!L38
$ws.Cells.Item(38, 13).Value = -1552.5  # ARM!M38
$ws.Cells.Item(38, 14).Value = -30934  # ARM!N38

$ws.Cells.Item(44, 8).Value = 28554.334  # ARM!H44
$ws.Cells.Item(44, 10).Value = 28554.334  # ARM!J44
$ws.Cells.Item(44, 12).Value = 28554.334  # ARM!L44
$ws.Cells.Item(44, 14).Value = -29530.334  # ARM!N44

$ws.Cells.Item(55, 8).Value = 19264.5  # ARM!H55
$ws.Cells.Item(55, 10).Value = 19264.5  # ARM!J55
$ws.Cells.Item(55, 12).Value = 19264.5  # ARM!L55
$ws.Cells.Item(55, 14).Value = -19894.5  # ARM!N55

$ws.Cells.Item(80, 8).Value = 34290  # ARM!H80
$ws.Cells.Item(80, 9).Value = 0  # ARM!I80
$ws.Cells.Item(80, 10).Value = 34290  # ARM!J80
$ws.Cells.Item(80, 11).Value = 0  # ARM!K80
$ws.Cells.Item(80, 12).Value = 34290  # ARM!L80
$ws.Cells.Item(80, 13).ClearContents()  # ARM!M80
$ws.Cells.Item(80, 14).Value = -36286  # ARM!N80

$ws.Cells.Item(83, 8).Value = 34290  # ARM!H83
$ws.Cells.Item(83, 9).Value = 0  # ARM!I83
$ws.Cells.Item(83, 10).Value = 34290  # ARM!J83
$ws.Cells.Item(83, 11).Value = 0  # ARM!K83
$ws.Cells.Item(83, 12).Value = 102870  # ARM!L83
$ws.Cells.Item(83, 13).ClearContents()  # ARM!M83
$ws.Cells.Item(83, 14).Value = -112854  # ARM!N83

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(33, 8).Value = 20298  # BSM!H33
$ws.Cells.Item(33, 9).Value = 25291.5  # BSM!I33
$ws.Cells.Item(33, 10).Value = 324  # BSM!J33
$ws.Cells.Item(33, 11).Value = 25291.5  # BSM!K33
$ws.Cells.Item(33, 12).Value = 324  # BSM!L33
$ws.Cells.Item(33, 13).Value = -24955.5  # BSM!M33
$ws.Cells.Item(33, 14).Value = -996  # BSM!N33

$ws.Cells.Item(134, 8).Value = 2625.348  # BSM!H134
$ws.Cells.Item(134, 9).Value = 2470.6191  # BSM!I134
$ws.Cells.Item(134, 10).Value = 4250  # BSM!J134
$ws.Cells.Item(134, 11).Value = 7411.8573  # BSM!K134
$ws.Cells.Item(134, 12).Value = 12750  # BSM!L134
$ws.Cells.Item(134, 13).Value = -4876.8573  # BSM!M134
$ws.Cells.Item(134, 14).Value = -17820  # BSM!N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 1782423.5  # CRP!H12
$ws.Cells.Item(12, 9).Value = 1669467.5  # CRP!I12
$ws.Cells.Item(12, 10).Value = 2008335.4  # CRP!J12
$ws.Cells.Item(12, 11).Value = 1669467.5  # CRP!K12
$ws.Cells.Item(12, 12).Value = 2008335.4  # CRP!L12
$ws.Cells.Item(12, 13).Value = -1669297.5  # CRP!M12
$ws.Cells.Item(12, 14).Value = -2008675.4  # CRP!N12

$ws.Cells.Item(14, 8).Value = 2000  # CRP!H14
$ws.Cells.Item(14, 10).Value = 2000  # CRP!J14
$ws.Cells.Item(14, 12).Value = 2000  # CRP!L14
$ws.Cells.Item(14, 14).Value = -2340  # CRP!N14

$ws.Cells.Item(32, 8).Value = 7505  # CRP!H32
$ws.Cells.Item(32, 9).Value = 7505  # CRP!I32
$ws.Cells.Item(32, 11).Value = 7505  # CRP!K32
$ws.Cells.Item(32, 13).Value = -7189  # CRP!M32

$ws.Cells.Item(35, 8).Value = 6120.769  # CRP!H35
$ws.Cells.Item(35, 9).Value = 1841.1111  # CRP!I35
$ws.Cells.Item(35, 10).Value = 15750  # CRP!J35
$ws.Cells.Item(35, 11).Value = 1841.1111  # CRP!K35
$ws.Cells.Item(35, 12).Value = 15750  # CRP!L35
$ws.Cells.Item(35, 13).Value = -1547.1111  # CRP!M35
$ws.Cells.Item(35, 14).Value = -16338  # CRP!N35

$ws.Cells.Item(38, 8).Value = 1940.5834  # CRP!H38
$ws.Cells.Item(38, 9).Value = 1220.4  # CRP!I38
$ws.Cells.Item(38, 10).Value = 5541.5  # CRP!J38
$ws.Cells.Item(38, 11).Value = 1220.4  # CRP!K38
$ws.Cells.Item(38, 12).Value = 5541.5  # CRP!L38
$ws.Cells.Item(38, 13).Value = -843.4000000000001  # CRP!M38
$ws.Cells.Item(38, 14).Value = -6295.5  # CRP!N38

$ws.Cells.Item(46, 8).Value = 1940.5834  # CRP!H46
$ws.Cells.Item(46, 9).Value = 1220.4  # CRP!I46
$ws.Cells.Item(46, 10).Value = 5541.5  # CRP!J46
$ws.Cells.Item(46, 11).Value = 1220.4  # CRP!K46
$ws.Cells.Item(46, 12).Value = 5541.5  # CRP!L46
$ws.Cells.Item(46, 13).Value = -1009.4  # CRP!M46
$ws.Cells.Item(46, 14).Value = -5963.5  # CRP!N46

$ws.Cells.Item(50, 8).Value = 13071  # CRP!H50
$ws.Cells.Item(50, 10).Value = 13071  # CRP!J50
$ws.Cells.Item(50, 12).Value = 13071  # CRP!L50
$ws.Cells.Item(50, 14).Value = -14321  # CRP!N50

$ws.Cells.Item(60, 8).Value = 19875.75  # CRP!H60
$ws.Cells.Item(60, 10).Value = 24501  # CRP!J60
$ws.Cells.Item(60, 12).Value = 24501  # CRP!L60
$ws.Cells.Item(60, 14).Value = -25523  # CRP!N60

$ws.Cells.Item(68, 8).Value = 27626.428  # CRP!H68
$ws.Cells.Item(68, 10).Value = 27626.428  # CRP!J68
$ws.Cells.Item(68, 12).Value = 27626.428  # CRP!L68
$ws.Cells.Item(68, 14).Value = -29124.428  # CRP!N68

$ws.Cells.Item(71, 8).Value = 27626.428  # CRP!H71
$ws.Cells.Item(71, 10).Value = 27626.428  # CRP!J71
$ws.Cells.Item(71, 12).Value = 82879.284  # CRP!L71
$ws.Cells.Item(71, 14).Value = -90367.284  # CRP!N71

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 75.42856999999999  # CUL!H8
$ws.Cells.Item(8, 9).Value = 75.42856999999999  # CUL!I8
$ws.Cells.Item(8, 11).Value = 226.28571  # CUL!K8
$ws.Cells.Item(8, 13).Value = -87.28570999999999  # CUL!M8

$ws.Cells.Item(16, 8).Value = 5632.4443  # CUL!H16
$ws.Cells.Item(16, 9).Value = 4599.75  # CUL!I16
$ws.Cells.Item(16, 10).Value = 6458.6  # CUL!J16
$ws.Cells.Item(16, 11).Value = 13799.25  # CUL!K16
$ws.Cells.Item(16, 12).Value = 19375.8  # CUL!L16
$ws.Cells.Item(16, 13).Value = -13626.25  # CUL!M16
$ws.Cells.Item(16, 14).Value = -19721.8  # CUL!N16

$ws.Cells.Item(20, 8).Value = 10197.333  # CUL!H20
$ws.Cells.Item(20, 10).Value = 10197.333  # CUL!J20
$ws.Cells.Item(20, 12).Value = 30591.999  # CUL!L20
$ws.Cells.Item(20, 14).Value = -31045.999  # CUL!N20

$ws.Cells.Item(31, 8).Value = 600  # CUL!H31
$ws.Cells.Item(31, 9).Value = 600  # CUL!I31
$ws.Cells.Item(31, 10).Value = 0  # CUL!J31
$ws.Cells.Item(31, 11).Value = 1800  # CUL!K31
$ws.Cells.Item(31, 12).Value = 0  # CUL!L31
$ws.Cells.Item(31, 13).Value = -1512  # CUL!M31
$ws.Cells.Item(31, 14).ClearContents()  # CUL!N31

$ws.Cells.Item(75, 8).Value = 15250  # CUL!H75
$ws.Cells.Item(75, 9).Value = 500  # CUL!I75
$ws.Cells.Item(75, 10).Value = 30000  # CUL!J75
$ws.Cells.Item(75, 11).Value = 1500  # CUL!K75
$ws.Cells.Item(75, 12).Value = 90000  # CUL!L75
$ws.Cells.Item(75, 13).Value = -502  # CUL!M75
$ws.Cells.Item(75, 14).Value = -91996  # CUL!N75

$ws.Cells.Item(78, 8).Value = 15250  # CUL!H78
$ws.Cells.Item(78, 9).Value = 500  # CUL!I78
$ws.Cells.Item(78, 10).Value = 30000  # CUL!J78
$ws.Cells.Item(78, 11).Value = 4500  # CUL!K78
$ws.Cells.Item(78, 12).Value = 270000  # CUL!L78
$ws.Cells.Item(78, 13).Value = 492  # CUL!M78
$ws.Cells.Item(78, 14).Value = -279984  # CUL!N78

$ws.Cells.Item(131, 8).Value = 1051.5454  # CUL!H131
$ws.Cells.Item(131, 10).Value = 1087.129  # CUL!J131
$ws.Cells.Item(131, 12).Value = 3261.387  # CUL!L131
$ws.Cells.Item(131, 14).Value = -13341.387  # CUL!N131

$ws.Cells.Item(138, 8).Value = 62501348  # CUL!H138
$ws.Cells.Item(138, 10).Value = 3000  # CUL!J138
$ws.Cells.Item(138, 12).Value = 9000  # CUL!L138
$ws.Cells.Item(138, 14).Value = -19280  # CUL!N138

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 2585926  # GSM!H24
$ws.Cells.Item(24, 10).Value = 3456.3333  # GSM!J24
$ws.Cells.Item(24, 12).Value = 3456.3333  # GSM!L24
$ws.Cells.Item(24, 14).Value = -3802.3333  # GSM!N24

$ws.Cells.Item(95, 8).Value = 23447.6  # GSM!H95
$ws.Cells.Item(95, 10).Value = 23447.6  # GSM!J95
$ws.Cells.Item(95, 12).Value = 23447.6  # GSM!L95
$ws.Cells.Item(95, 14).Value = -28939.6  # GSM!N95

$ws.Cells.Item(122, 8).Value = 5682.1  # GSM!H122
$ws.Cells.Item(122, 9).Value = 4727.625  # GSM!I122
$ws.Cells.Item(122, 10).Value = 9500  # GSM!J122
$ws.Cells.Item(122, 11).Value = 14182.875  # GSM!K122
$ws.Cells.Item(122, 12).Value = 28500  # GSM!L122
$ws.Cells.Item(122, 13).Value = -11732.875  # GSM!M122
$ws.Cells.Item(122, 14).Value = -33400  # GSM!N122

$ws.Cells.Item(123, 8).Value = 9993.333000000001  # GSM!H123
$ws.Cells.Item(123, 10).Value = 9993.333000000001  # GSM!J123
$ws.Cells.Item(123, 12).Value = 9993.333000000001  # GSM!L123
$ws.Cells.Item(123, 14).Value = -14893.333  # GSM!N123

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2362.75  # LTW!H7
$ws.Cells.Item(7, 9).Value = 1272  # LTW!I7
$ws.Cells.Item(7, 10).Value = 3889.8  # LTW!J7
$ws.Cells.Item(7, 11).Value = 1272  # LTW!K7
$ws.Cells.Item(7, 12).Value = 3889.8  # LTW!L7
$ws.Cells.Item(7, 13).Value = -1160  # LTW!M7
$ws.Cells.Item(7, 14).Value = -4113.8  # LTW!N7

$ws.Cells.Item(22, 8).Value = 738.0454999999999  # LTW!H22
$ws.Cells.Item(22, 9).Value = 566.0833  # LTW!I22
$ws.Cells.Item(22, 11).Value = 566.0833  # LTW!K22
$ws.Cells.Item(22, 13).Value = -271.0833  # LTW!M22

$ws.Cells.Item(27, 8).Value = 738.0454999999999  # LTW!H27
$ws.Cells.Item(27, 9).Value = 566.0833  # LTW!I27
$ws.Cells.Item(27, 11).Value = 566.0833  # LTW!K27
$ws.Cells.Item(27, 13).Value = -459.0833  # LTW!M27

$ws.Cells.Item(32, 8).Value = 3231.5  # LTW!H32
$ws.Cells.Item(32, 9).Value = 2807.4  # LTW!I32
$ws.Cells.Item(32, 10).Value = 3938.3333  # LTW!J32
$ws.Cells.Item(32, 11).Value = 2807.4  # LTW!K32
$ws.Cells.Item(32, 12).Value = 3938.3333  # LTW!L32
$ws.Cells.Item(32, 13).Value = -2490.4  # LTW!M32
$ws.Cells.Item(32, 14).Value = -4572.3333  # LTW!N32

$ws.Cells.Item(126, 8).Value = 2362.75  # LTW!H126
$ws.Cells.Item(126, 9).Value = 1272  # LTW!I126
$ws.Cells.Item(126, 10).Value = 3889.8  # LTW!J126
$ws.Cells.Item(126, 11).Value = 3816  # LTW!K126
$ws.Cells.Item(126, 12).Value = 11669.4  # LTW!L126
$ws.Cells.Item(126, 13).Value = -1346  # LTW!M126
$ws.Cells.Item(126, 14).Value = -16609.4  # LTW!N126

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 0  # WVR!H109
$ws.Cells.Item(109, 10).Value = 0  # WVR!J109
$ws.Cells.Item(109, 12).Value = 0  # WVR!L109
$ws.Cells.Item(109, 14).ClearContents()  # WVR!N109

$ws.Cells.Item(132, 8).Value = 2500.9836  # WVR!H132
$ws.Cells.Item(132, 9).Value = 2563.878  # WVR!I132
$ws.Cells.Item(132, 11).Value = 7691.634  # WVR!K132
$ws.Cells.Item(132, 13).Value = -5161.634  # WVR!M132
